$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute(
    "Exploring the Wonders of the Quantum Realm", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Embracing Diversity in a Globalized World: Navigating Challenges and Celebrating Unity", 2)

# --- Author name (collapses the 5 runs "Dr" "." " Albert J" "." " Robertson" into one) ---
$d.Content.Find.Execute(
    "Dr. Albert J. Robertson", $false, $false, $false, $false, $false,
    $true, 1, $false, "Hannah Davies", 2)

# --- Email address local-part / domain (keep the existing ".org") ---
$d.Content.Find.Execute(
    "quantumstudies@scientificdiscovery", $false, $false, $false, $false, $false,
    $true, 1, $false, "hannah.davies@eduworld", 2)

# --- Body paragraph (paragraph 5): full rewrite, preserving embedded line breaks ---
$body = "In a world where borders are increasingly blurred and cultures intertwine, embracing diversity has become a cornerstone of global harmony. As citizens of a globalized society, it is imperative that we understand, appreciate, and celebrate the rich tapestry of cultures, beliefs, and perspectives that make our world a vibrant and dynamic place. In this essay, we will delve into the multifaceted nature of diversity, explore the challenges it presents, and highlight the immense benefits that stem from fostering inclusivity and understanding." + [char]11 + "" + [char]11 + "The beauty of diversity lies in its multifaceted nature. It encompasses differences in race, ethnicity, gender, sexual orientation, religion, language, and socioeconomic status, among others. Each individual brings a unique set of experiences, perspectives, and talents to the table, creating a rich and dynamic social fabric. Diversity challenges us to step outside our comfort zones, question our assumptions, and engage with those who hold different beliefs and values. By doing so, we broaden our horizons, foster empathy, and gain a deeper understanding of the human experience." + [char]11 + "" + [char]11 + "However, embracing diversity is not without its challenges. Misunderstandings, prejudice, and discrimination can arise when people from different backgrounds interact. These challenges can be daunting, but they also present opportunities for growth and learning. By promoting open dialogue, encouraging inclusivity, and challenging stereotypes, we can create a society where diversity is celebrated and everyone feels valued and respected." + [char]11 + "" + [char]11 + "Introduction Continued:" + [char]11 + "" + [char]11 + "The immense benefits that stem from fostering inclusivity and understanding are undeniable. A diverse society is a more vibrant, innovative, and resilient one. When people from different backgrounds come together, they bring fresh ideas, perspectives, and solutions to complex problems. This cross-pollination of ideas leads to groundbreaking discoveries, transformative technologies, and creative expressions that benefit all of humanity. Furthermore, a diverse society is a more tolerant and peaceful one. When people understand and respect each other's differences, they are less likely to resort to violence or conflict. Instead, they are more likely to work together to build a better future for themselves and for generations to come." + [char]11 + "" + [char]11 + "Introduction Continued:" + [char]11 + "" + [char]11 + "The journey towards a truly diverse and inclusive society is an ongoing one, but it is a journey worth taking. By challenging our biases, educating ourselves about different cultures, and actively promoting inclusivity, we can create a world where everyone feels valued, respected, and empowered. A world where diversity is not just tolerated but celebrated, and where the unique contributions of each individual are recognized and appreciated. In this world, the challenges of diversity will be outweighed by the boundless opportunities it presents, leading to a more harmonious, equitable, and prosperous global community."
$bodyRange = $d.Paragraphs(5).Range
$bodyRange.End = $bodyRange.End - 1
$bodyRange.Text = $body

# --- Summary paragraph (paragraph 7): full rewrite ---
$summary = "In conclusion, embracing diversity in a globalized world presents both challenges and immense benefits. By fostering inclusivity, understanding, and respect, we can create a society where everyone feels valued and empowered. The challenges of diversity, such as misunderstandings and prejudice, can be overcome through open dialogue, education, and a commitment to building a more just and equitable world. The beauty of diversity lies in its multifaceted nature, and it is through celebrating our differences and coming together as a global community that we can unlock the full potential of our shared humanity."
$summaryRange = $d.Paragraphs(7).Range
$summaryRange.End = $summaryRange.End - 1
$summaryRange.Text = $summary

# --- Append a new empty paragraph at the end of the document ---
$d.Paragraphs(7).Range.InsertParagraphAfter()

Write-Output "done"
